# Add LEDs for Power and status indication (LED1, LED2), and a new resistor
# line (R28, R30) that shares the footprint/value family with R4/R12 (which
# picks up two more designators: R27, R29).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Insert two new rows at row 20 (just above "Q1, Q2, Q3") for the new
#    LED1 / LED2 BOM lines. Insert one at a time so each inherits the
#    formatting of the row immediately above it (row 19 = "J9").
# ---------------------------------------------------------------------
$ws.Rows.Item(20).Insert()
$ws.Rows.Item(21).Insert()

# LED1
$ws.Range("A20").Value = "LED1"
$ws.Range("B20").Value = 1
$ws.Range("C20").Value = "20mA 54mcd -55?~+85? 631nm Red 130° 75mW 2V 0603 Light Emitting Diodes (LED) ROHS"
$ws.Range("D20").Value = "Lite-On"
$ws.Range("E20").Value = "LTST-C191KRKT"
$ws.Range("F20").Value = "'0.014"
$ws.Range("G20").Value = "'0.011"
$ws.Range("H20").Value = "LCSC"
$ws.Range("I20").Value = "C125099"

# LED2
$ws.Range("A21").Value = "LED2"
$ws.Range("B21").Value = 1
$ws.Range("C21").Value = "SMD5050 RGB LEDs ROHS"
$ws.Range("D21").Value = "HONGLITRONIC"
$ws.Range("E21").Value = "HL-AF-5060H248BS36FU76GC-S1-THL"
$ws.Range("F21").Value = "'0.014"
$ws.Range("G21").Value = "'0.0325"
$ws.Range("H21").Value = "LCSC"
$ws.Range("I21").Value = "C2683773"

# Re-apply the formatting (border/shading/quote-prefix style) from the row
# directly above the inserted block so the new rows visually match the rest
# of the table, without disturbing the values we just wrote.
$ws.Range("A19:I19").Copy()
$ws.Range("A20:I21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2) "R4, R12" (now shifted down to row 25) also covers two newly added
#    resistors, R27 and R29 -> designator + quantity update only.
# ---------------------------------------------------------------------
$ws.Range("A25").Value = "R4, R12, R27, R29"
$ws.Range("B25").Value = 4

# ---------------------------------------------------------------------
# 3) Insert a new row after "R13" (now row 31) for "R28, R30".
# ---------------------------------------------------------------------
$ws.Rows.Item(32).Insert()

$ws.Range("A32").Value = "R28, R30"
$ws.Range("B32").Value = 2
$ws.Range("C32").Value = "100mW Thick Film Resistors 75V ±100ppm/? ±1% 30O 0603 Chip Resistor - Surface Mount ROHS"
$ws.Range("D32").Value = "Yageo"
$ws.Range("E32").Value = "RC0603FR-0730RL"
$ws.Range("F32").Value = "'0.0012"
$ws.Range("G32").Value = "'0.001"
$ws.Range("H32").Value = "LCSC"
$ws.Range("I32").Value = "C128060"

$ws.Range("A31:I31").Copy()
$ws.Range("A32:I32").PasteSpecial(-4122)
$excel.CutCopyMode = $false

"done"
